$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Sample_name" column at B (shifts existing B..K one column right to C..L)
$ws.Columns("B").Insert()

# Header
$ws.Range("B1").Value = "Sample_name"

# Sample names per row (aligned with each row's CELL_DIVE_ID / Response)
$ws.Range("B2").Value = "PR"
$ws.Range("B3").Value = "CR"
$ws.Range("B4").Value = "Untreated"

# New column reads about as wide as the CELL_DIVE_ID column (A)
$ws.Columns("B").ColumnWidth = 19.6
